$d = $word.ActiveDocument

$d.Content.Find.Execute("217÷7=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "220÷4=55, 0", 2) | Out-Null
$d.Content.Find.Execute("236÷9=26, 2", $true, $false, $false, $false, $false, $true, 1, $false, "771÷3=257, 0", 2) | Out-Null
$d.Content.Find.Execute("839÷4=209, 3", $true, $false, $false, $false, $false, $true, 1, $false, "332÷4=83, 0", 2) | Out-Null
$d.Content.Find.Execute("838÷6=139, 4", $true, $false, $false, $false, $false, $true, 1, $false, "331÷5=66, 1", 2) | Out-Null
$d.Content.Find.Execute("299÷5=59, 4", $true, $false, $false, $false, $false, $true, 1, $false, "391÷6=65, 1", 2) | Out-Null
$d.Content.Find.Execute("632÷8=79, 0", $true, $false, $false, $false, $false, $true, 1, $false, "236÷4=59, 0", 2) | Out-Null
$d.Content.Find.Execute("540÷3=180, 0", $true, $false, $false, $false, $false, $true, 1, $false, "247÷3=82, 1", 2) | Out-Null
$d.Content.Find.Execute("103÷7=14, 5", $true, $false, $false, $false, $false, $true, 1, $false, "533÷6=88, 5", 2) | Out-Null
$d.Content.Find.Execute("467÷7=66, 5", $true, $false, $false, $false, $false, $true, 1, $false, "527÷8=65, 7", 2) | Out-Null
$d.Content.Find.Execute("170÷5=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "525÷9=58, 3", 2) | Out-Null
$d.Content.Find.Execute("830÷9=92, 2", $true, $false, $false, $false, $false, $true, 1, $false, "194÷6=32, 2", 2) | Out-Null
$d.Content.Find.Execute("845÷3=281, 2", $true, $false, $false, $false, $false, $true, 1, $false, "992÷9=110, 2", 2) | Out-Null
$d.Content.Find.Execute("410÷9=45, 5", $true, $false, $false, $false, $false, $true, 1, $false, "737÷6=122, 5", 2) | Out-Null
$d.Content.Find.Execute("575÷2=287, 1", $true, $false, $false, $false, $false, $true, 1, $false, "978÷3=326, 0", 2) | Out-Null
$d.Content.Find.Execute("352÷2=176, 0", $true, $false, $false, $false, $false, $true, 1, $false, "907÷6=151, 1", 2) | Out-Null
$d.Content.Find.Execute("753÷9=83, 6", $true, $false, $false, $false, $false, $true, 1, $false, "758÷4=189, 2", 2) | Out-Null
$d.Content.Find.Execute("417÷2=208, 1", $true, $false, $false, $false, $false, $true, 1, $false, "262÷8=32, 6", 2) | Out-Null
$d.Content.Find.Execute("511÷9=56, 7", $true, $false, $false, $false, $false, $true, 1, $false, "442÷7=63, 1", 2) | Out-Null
$d.Content.Find.Execute("767÷9=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "144÷5=28, 4", 2) | Out-Null
$d.Content.Find.Execute("318÷9=35, 3", $true, $false, $false, $false, $false, $true, 1, $false, "513÷4=128, 1", 2) | Out-Null
$d.Content.Find.Execute("820÷6=136, 4", $true, $false, $false, $false, $false, $true, 1, $false, "543÷5=108, 3", 2) | Out-Null
$d.Content.Find.Execute("803÷7=114, 5", $true, $false, $false, $false, $false, $true, 1, $false, "229÷6=38, 1", 2) | Out-Null
$d.Content.Find.Execute("916÷9=101, 7", $true, $false, $false, $false, $false, $true, 1, $false, "359÷3=119, 2", 2) | Out-Null
$d.Content.Find.Execute("322÷2=161, 0", $true, $false, $false, $false, $false, $true, 1, $false, "885÷6=147, 3", 2) | Out-Null
$d.Content.Find.Execute("368÷6=61, 2", $true, $false, $false, $false, $false, $true, 1, $false, "578÷7=82, 4", 2) | Out-Null
